{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items,text\");\nawait context.sync();\n\n// --- helpers --------------------------------------------------------\n\n// Build a minimal single-part OOXML payload usable with Range.insertOoxml /\n// Paragraph.insertOoxml. Using insertOoxml (instead of insertText) keeps the\n// inserted text in its own <w:r> run instead of letting the host merge it\n// into an adjacent run that happens to share the same formatting.\nfunction ooxmlRun(text) {\n  const esc = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  const needsPreserve = /^\\s|\\s$/.test(text);\n  const spaceAttr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t${spaceAttr}>${esc}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n}\n\n// Appends `text` to the end of paragraph `p` as a brand-new run.\nasync function appendRun(p, text) {\n  const tail = p.getRange(\"End\");\n  tail.insertOoxml(ooxmlRun(text), Word.InsertLocation.end);\n  await context.sync();\n}\n\n// --- original paragraphs (0-indexed) --------------------------------\n// 0: \"Binance passes:\"\n// 1: \"\" (empty)\n// 2: \"Pass: +Jamilchik1975\"\n// 3: \"Passkey: 608053\"\n// 4: \"Wallet1: +Jamilchik1975\"\n\nconst p0 = paras.items[0];\nconst p1 = paras.items[1];\nconst p2 = paras.items[2];\nconst p3 = paras.items[3];\nconst p4 = paras.items[4];\n\n// 1) \"Binance passes:\" -> \"Bin passes:\"\np0.getRange(\"Whole\").insertText(\"Bin passes:\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Drop the blank paragraph that used to separate the heading from the\n//    \"Pass:\" line \u2014 everything below gets folded into paragraph 0 instead.\np1.delete();\nawait context.sync();\n\n// 3) Rebuild paragraph 0 as: \"Bin passes:\" \" \" \"Pass: +Jamilchik1975\" \" / \"\n//    \"Passkey: 608053\" \" / \" \"Wallet1: +Jamilchik1975\" (each segment its own run).\nawait appendRun(p0, \" \");\nawait appendRun(p0, \"Pass: +Jamilchik1975\");\nawait appendRun(p0, \" / \");\nawait appendRun(p0, \"Passkey: 608053\");\nawait appendRun(p0, \" / \");\nawait appendRun(p0, \"Wallet1: +Jamilchik1975\");\n\n// 4) The standalone \"Pass: +Jamilchik1975\" paragraph is now redundant\n//    (its text lives inline in paragraph 0) \u2014 remove it.\np2.delete();\nawait context.sync();\n\n// 5) Recycle the old \"Passkey: 608053\" paragraph to read \"Gmail: ...\"\np3.getRange(\"Whole\").insertText(\"Gmail: Jamilchik1975\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 6) Recycle the old \"Wallet1: +Jamilchik1975\" paragraph to read \"Hotmail: ...\"\np4.getRange(\"Whole\").insertText(\n  \"Hotmail: Jamil1975 or +Jamil1975 or Jamil0123\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 7) Append the rest of the new credential lines after that paragraph.\nconst facebook = p4.insertParagraph(\"Facebook: 21021975\", Word.InsertLocation.after);\nawait context.sync();\n\nconst insta = facebook.insertParagraph(\"Insta: jamilchik1111\", Word.InsertLocation.after);\nawait context.sync();\n\nconst deco = insta.insertParagraph(\"DECO: Jamilchik1975\", Word.InsertLocation.after);\nawait context.sync();\n\nconst home123 = deco.insertParagraph(\"Home123: 12345678\", Word.InsertLocation.after);\nawait context.sync();\n\n// \"GAKhome\" is flagged by Word's proofing engine as a misspelling, so the\n// authored document wraps it in spell-check markers; reproduce that here.\nconst gakhome = home123.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\nconst gakhomeOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>GAKhome</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t>: 04522966</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\ngakhome.getRange(\"Whole\").insertOoxml(gakhomeOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\nconst laptop = gakhome.insertParagraph(\"Laptop: 0000\", Word.InsertLocation.after);\nawait context.sync();\n\nreturn \"done\";\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is pre-seeded as $d below.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Helper: build a minimal single-part WordOpenXML \"package\" payload that\n# InsertXML understands, containing one <w:p> with one <w:r> per\n# supplied string. Using InsertXML (instead of Range.Text / InsertAfter)\n# keeps each piece of text in its own <w:r> run instead of letting\n# same-formatted adjacent text coalesce into a single run.\n# ---------------------------------------------------------------------\nfunction New-ParaXml([string[]]$texts) {\n    $runs = \"\"\n    foreach ($t in $texts) {\n        $esc = $t -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n        if ($t -match '^\\s' -or $t -match '\\s$') {\n            $runs += '<w:r><w:t xml:space=\"preserve\">' + $esc + '</w:t></w:r>'\n        } else {\n            $runs += '<w:r><w:t>' + $esc + '</w:t></w:r>'\n        }\n    }\n    return '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runs + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Original paragraphs (1-indexed, Word COM convention):\n# 1: \"Binance passes:\"\n# 2: \"\" (empty)\n# 3: \"Pass: +Jamilchik1975\"\n# 4: \"Passkey: 608053\"\n# 5: \"Wallet1: +Jamilchik1975\"\n\n# 1) Rebuild paragraph 1's text-only range (Start..End-1, excluding the\n#    trailing paragraph mark) as multiple runs: \"Bin passes:\" \" \"\n#    \"Pass: +Jamilchik1975\" \" / \" \"Passkey: 608053\" \" / \"\n#    \"Wallet1: +Jamilchik1975\".\n$p1 = $d.Paragraphs.Item(1)\n$textRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)\n$textRange.InsertXML((New-ParaXml @(\n    \"Bin passes:\",\n    \" \",\n    \"Pass: +Jamilchik1975\",\n    \" / \",\n    \"Passkey: 608053\",\n    \" / \",\n    \"Wallet1: +Jamilchik1975\"\n)))\n\n# 2) Remove the (now second) empty paragraph that used to separate the\n#    heading from the \"Pass:\" line.\n$emptyPara = $d.Paragraphs.Item(2)\n$emptyPara.Range.Delete()\n\n# 3) The old standalone \"Pass: +Jamilchik1975\" paragraph (now item 2) is\n#    redundant since its text lives inline in paragraph 1 \u2014 remove it.\n$oldPass = $d.Paragraphs.Item(2)\n$oldPass.Range.Delete()\n\n# 4) Recycle the old \"Passkey: 608053\" paragraph (now item 2) to read\n#    \"Gmail: Jamilchik1975\".\n$gmailPara = $d.Paragraphs.Item(2)\n$gmailRange = $d.Range($gmailPara.Range.Start, $gmailPara.Range.End - 1)\n$gmailRange.InsertXML((New-ParaXml @(\"Gmail: Jamilchik1975\")))\n\n# 5) Recycle the old \"Wallet1: +Jamilchik1975\" paragraph (now item 3) to\n#    read the Hotmail line.\n$hotmailPara = $d.Paragraphs.Item(3)\n$hotmailRange = $d.Range($hotmailPara.Range.Start, $hotmailPara.Range.End - 1)\n$hotmailRange.InsertXML((New-ParaXml @(\"Hotmail: Jamil1975 or +Jamil1975 or Jamil0123\")))\n\n# 6) Append the remaining new credential paragraphs at the end of the\n#    document.\n$endRange = $d.Content\n$endRange.Collapse(0)  # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange = $d.Paragraphs.Last.Range\n$endRange.Text = \"Facebook: 21021975\"\n\n$p = $d.Paragraphs.Last.Range\n$p.InsertParagraphAfter()\n$p = $d.Paragraphs.Last.Range\n$p.Text = \"Insta: jamilchik1111\"\n\n$p = $d.Paragraphs.Last.Range\n$p.InsertParagraphAfter()\n$p = $d.Paragraphs.Last.Range\n$p.Text = \"DECO: Jamilchik1975\"\n\n$p = $d.Paragraphs.Last.Range\n$p.InsertParagraphAfter()\n$p = $d.Paragraphs.Last.Range\n$p.Text = \"Home123: 12345678\"\n\n$p = $d.Paragraphs.Last.Range\n$p.InsertParagraphAfter()\n$p = $d.Paragraphs.Last.Range\n# Seed placeholder text so the range is non-empty, then splice in the real\n# runs via InsertXML over that text range (mirrors the technique used\n# above \u2014 InsertXML on a collapsed/empty range inserts an extra paragraph\n# break instead of replacing in place).\n$p.Text = \"GAKhome: 04522966\"\n$gakPara = $d.Paragraphs.Last\n$gakRange = $d.Range($gakPara.Range.Start, $gakPara.Range.End - 1)\n# \"GAKhome\" is flagged by Word's proofing engine as a misspelling, so the\n# authored document wraps it in spell-check markers; reproduce that here.\n$gakXml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>GAKhome</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>: 04522966</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$gakRange.InsertXML($gakXml)\n\n$p = $d.Paragraphs.Last.Range\n$p.InsertParagraphAfter()\n$p = $d.Paragraphs.Last.Range\n$p.Text = \"Laptop: 0000\"\n"}
